# "rectified date picker issues"
#
# The V2Project sheet's row 2 holds a single sample "project" record used to
# exercise the date-picker / registration flow in the front-end tests. The
# record's identifying fields were stale, so refresh them with a newly
# generated test project (new ProjectID, ProjectName, State, Street, City,
# PostalCode and Area) while leaving every other column (Country, Website,
# contact info, ...) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2Project")

# ProjectID (A2) is a purely-numeric string; prefix with an apostrophe so it
# is stored as text (matching the sheet's existing quote-prefixed text
# style) instead of being coerced to a number.
$ws.Range("A2").Value = "'2202278615"

# ProjectName (B2)
$ws.Range("B2").Value = "Automation V2 Project7288747"

# State (E2)
$ws.Range("E2").Value = "South Carolina"

# Street (F2)
$ws.Range("F2").Value = "13210 Young Crossroad"

# City (G2)
$ws.Range("G2").Value = "Carterside"

# PostalCode (H2) - also purely numeric; force text formatting so it is
# stored as a string (not a number), then drop back to the sheet's plain
# (unstyled) look so no stray cell formatting is introduced.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "74228"
$ws.Range("H2").Style = "Normal"

# Area (J2) - same purely-numeric-text situation as PostalCode.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "46899"
$ws.Range("J2").Style = "Normal"
